$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: add a new (blank) "Paid" column --------
# A new column is inserted before column N, shifting the existing
# "Late" / heading / "Outstanding" columns one place to the right
# (N->O, O->P, P->Q). The new column inherits the width of the column
# immediately to its left (column M).
$ws = $wb.Worksheets.Item("Repayment schedule")

$leftColumnWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftColumnWidth

# --- Make "Repayment schedule" the active sheet/tab ----------------------
$ws.Activate()
$ws.Range("R11").Select()
